$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the dSF column (F) values that were repulled/recalculated.
$ws.Range("F2").Value = 5
$ws.Range("F5").Value = 4
$ws.Range("F12").Value = -8
$ws.Range("F13").Value = -2
